# Trade #64 closed at 2026-02-17 12:53:28 - unknown UNKNOWN +0.000%
#
# - Updates the aggregate metrics on the "Summary" sheet and on the
#   "MarketMaking" row of the "Strategy Status" sheet to reflect the
#   newly closed trade.
# - Appends the new trade (#64) as row 65 on both the "All Trades" and
#   "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.11   # Current Capital
$summary.Range("B4").Value = 0.1       # Total P&L $
$summary.Range("B6").Value = 64        # Total Trades
$summary.Range("B8").Value = 21        # Losing Trades
$summary.Range("B9").Value = 43.75     # Win Rate %

# --- Strategy Status sheet (MarketMaking is row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.11     # Capital
$status.Range("D4").Value = 64         # Trades
$status.Range("E4").Value = 0.1        # P&L $
$status.Range("F4").Value = 0.11       # P&L %
$status.Range("G4").Value = 43.75      # Win Rate %

# --- Append new trade row (row 65) to "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A65").Value = 64

    # The Date column looks like "2026-02-17" which Excel would otherwise
    # auto-detect and silently convert into a date serial number. Force it
    # to be entered as plain text, matching the source data, then restore
    # the cell's default (Normal) style so no stray number format lingers.
    $ws.Range("B65").NumberFormat = "@"
    $ws.Range("B65").Value = "2026-02-17"
    $ws.Range("B65").Style = "Normal"

    $ws.Range("C65").Value = "12:53:21"
    $ws.Range("D65").Value = "MarketMaking"
    $ws.Range("E65").Value = "DOWN"
    $ws.Range("F65").Value = 0.88
    $ws.Range("G65").Value = 0.87
    $ws.Range("H65").Value = "CLOSED"
    $ws.Range("I65").Value = -1.1364
    $ws.Range("J65").Value = -0.01
    $ws.Range("K65").Value = 100.11
    $ws.Range("L65").Value = 0
    $ws.Range("M65").Value = 0
    $ws.Range("N65").Value = 0.6
    $ws.Range("O65").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P65").Value = "early_exit"
    $ws.Range("Q65").Value = 0.13
}
